$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E33").Value = "N=3 n=161; Random-effects: 28.8% vs. 4.1%, OR=9.55, 95%CI: 3.5, 26.09; Fixed-effecs: OR=9.55, 95%CI: 3.5, 26.09"
$ws.Range("F33").Value = "66.7% had overall low risk of bias"
$ws.Range("G33").Value = "25% of eligible studies and 25.7% of participants had usable data"
$ws.Range("H33").Value = "33.3% studies and 42.2% participants with schizophrenia"

$ws.Range("D34").Value = "4-6 weeks"
$ws.Range("E34").Value = "N=3 n=497; Random-effects: 5.6% vs. 6%, OR=0.93, 95%CI: 0.17, 5.06; Fixed-effecs: OR=0.71, 95%CI: 0.32, 1.55"
$ws.Range("G34").Value = "50% of eligible studies and 32.1% of participants had usable data"
$ws.Range("H34").Value = "66.7% studies and 92.2% participants with schizophrenia"
